$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, 5, 1),
  @(2, 6, 0.3333333333333333),
  @(2, 7, 0.1050873333333333),
  @(2, 8, 0.315262),
  @(2, 9, 0.03031434174852429),
  @(2, 10, 0.03031434174852429),
  @(2, 13, 1.716657),
  @(2, 14, 5.149971),
  @(2, 15, 0.3840886036988016),
  @(2, 16, 0.3840886036988015),
  @(2, 17, 0.180398906378),
  @(2, 18, 1.623590157402),
  @(2, 19, 0.01164339319423898),
  @(2, 20, 0.01164339319423898),
  @(3, 5, 1),
  @(3, 6, 0.3333333333333333),
  @(3, 7, 0.1050873333333333),
  @(3, 8, 0.315262),
  @(3, 9, 0.03031434174852429),
  @(3, 10, 0.03031434174852429),
  @(3, 15, 0.07870146593648156),
  @(3, 16, 0.07870146593648154),
  @(3, 17, 0.03696453955822222),
  @(3, 18, 0.332680856024),
  @(3, 19, 0.002385783134508345),
  @(3, 20, 0.002385783134508345),
  @(4, 5, 1),
  @(4, 6, 0.3333333333333333),
  @(4, 7, 0.1050873333333333),
  @(4, 8, 0.315262),
  @(4, 9, 0.03031434174852429),
  @(4, 10, 0.03031434174852429),
  @(4, 13, 1.677572333333333),
  @(4, 14, 5.032717),
  @(4, 15, 0.3753437146230962),
  @(4, 16, 0.3753437146230962),
  @(4, 17, 0.1762916029837777),
  @(4, 18, 1.586624426854),
  @(4, 19, 0.01137829763824511),
  @(4, 20, 0.01137829763824511),
  @(5, 5, 1),
  @(5, 6, 0.3333333333333333),
  @(5, 7, 0.1050873333333333),
  @(5, 8, 0.315262),
  @(5, 9, 0.03031434174852429),
  @(5, 10, 0.03031434174852429),
  @(5, 13, 0.7234496666666667),
  @(5, 14, 2.170349),
  @(5, 15, 0.1618662157416207),
  @(5, 16, 0.1618662157416207),
  @(5, 17, 0.07602539627088888),
  @(5, 18, 0.6842285664379999),
  @(5, 19, 0.004906867781531854),
  @(5, 20, 0.004906867781531853),
  @(6, 9, 0.6354599969768544),
  @(6, 10, 0.6354599969768545),
  @(6, 13, 1.716657),
  @(6, 14, 5.149971),
  @(6, 15, 0.3840886036988016),
  @(6, 16, 0.3840886036988015),
  @(6, 17, 3.781585938845999),
  @(6, 18, 34.034273449614),
  @(6, 19, 0.2440729429452847),
  @(6, 20, 0.2440729429452847),
  @(7, 9, 0.6354599969768544),
  @(7, 10, 0.6354599969768545),
  @(7, 15, 0.07870146593648156),
  @(7, 16, 0.07870146593648154),
  @(7, 19, 0.05001163330607058),
  @(7, 20, 0.05001163330607058),
  @(8, 9, 0.6354599969768544),
  @(8, 10, 0.6354599969768545),
  @(8, 13, 1.677572333333333),
  @(8, 14, 5.032717),
  @(8, 15, 0.3753437146230962),
  @(8, 16, 0.3753437146230962),
  @(8, 17, 3.695487186508666),
  @(8, 18, 33.25938467857799),
  @(8, 19, 0.238515915759674),
  @(8, 20, 0.238515915759674),
  @(9, 9, 0.6354599969768544),
  @(9, 10, 0.6354599969768545),
  @(9, 13, 0.7234496666666667),
  @(9, 14, 2.170349),
  @(9, 15, 0.1618662157416207),
  @(9, 16, 0.1618662157416207),
  @(9, 17, 1.593671354807333),
  @(9, 18, 14.343042193266),
  @(9, 19, 0.1028595049658252),
  @(9, 20, 0.1028595049658252),
  @(10, 5, 3),
  @(10, 6, 1),
  @(10, 7, 0.9666886666666668),
  @(10, 8, 2.900066),
  @(10, 9, 0.2788588279503266),
  @(10, 10, 0.2788588279503266),
  @(10, 13, 1.716657),
  @(10, 14, 5.149971),
  @(10, 15, 0.3840886036988016),
  @(10, 16, 0.3840886036988015),
  @(10, 17, 1.659472866454),
  @(10, 18, 14.935255798086),
  @(10, 19, 0.1071064978565253),
  @(10, 20, 0.1071064978565253),
  @(11, 5, 3),
  @(11, 6, 1),
  @(11, 7, 0.9666886666666668),
  @(11, 8, 2.900066),
  @(11, 9, 0.2788588279503266),
  @(11, 10, 0.2788588279503266),
  @(11, 15, 0.07870146593648156),
  @(11, 16, 0.07870146593648154),
  @(11, 17, 0.3400333829591112),
  @(11, 18, 3.060300446632001),
  @(11, 19, 0.0219465985490198),
  @(11, 20, 0.0219465985490198),
  @(12, 5, 3),
  @(12, 6, 1),
  @(12, 7, 0.9666886666666668),
  @(12, 8, 2.900066),
  @(12, 9, 0.2788588279503266),
  @(12, 10, 0.2788588279503266),
  @(12, 13, 1.677572333333333),
  @(12, 14, 5.032717),
  @(12, 15, 0.3753437146230962),
  @(12, 16, 0.3753437146230962),
  @(12, 17, 1.621690162146889),
  @(12, 18, 14.595211459322),
  @(12, 19, 0.1046679083383185),
  @(12, 20, 0.1046679083383185),
  @(13, 5, 3),
  @(13, 6, 1),
  @(13, 7, 0.9666886666666668),
  @(13, 8, 2.900066),
  @(13, 9, 0.2788588279503266),
  @(13, 10, 0.2788588279503266),
  @(13, 13, 0.7234496666666667),
  @(13, 14, 2.170349),
  @(13, 15, 0.1618662157416207),
  @(13, 16, 0.1618662157416207),
  @(13, 17, 0.6993505936704445),
  @(13, 18, 6.294155343034),
  @(13, 19, 0.04513782320646307),
  @(13, 20, 0.04513782320646306),
  @(14, 7, 0.191934),
  @(14, 8, 0.575802),
  @(14, 9, 0.05536683332429468),
  @(14, 10, 0.05536683332429467),
  @(14, 13, 1.716657),
  @(14, 14, 5.149971),
  @(14, 15, 0.3840886036988016),
  @(14, 16, 0.3840886036988015),
  @(14, 17, 0.329484844638),
  @(14, 18, 2.965363601742),
  @(14, 19, 0.02126576970275262),
  @(14, 20, 0.02126576970275261),
  @(15, 7, 0.191934),
  @(15, 8, 0.575802),
  @(15, 9, 0.05536683332429468),
  @(15, 10, 0.05536683332429467),
  @(15, 15, 0.07870146593648156),
  @(15, 16, 0.07870146593648154),
  @(15, 17, 0.06751291245600001),
  @(15, 18, 0.6076162121040001),
  @(15, 19, 0.00435745094688283),
  @(15, 20, 0.004357450946882828),
  @(16, 7, 0.191934),
  @(16, 8, 0.575802),
  @(16, 9, 0.05536683332429468),
  @(16, 10, 0.05536683332429467),
  @(16, 13, 1.677572333333333),
  @(16, 14, 5.032717),
  @(16, 15, 0.3753437146230962),
  @(16, 16, 0.3753437146230962),
  @(16, 17, 0.321983168226),
  @(16, 18, 2.897848514034),
  @(16, 19, 0.0207815928868586),
  @(16, 20, 0.02078159288685859),
  @(17, 7, 0.191934),
  @(17, 8, 0.575802),
  @(17, 9, 0.05536683332429468),
  @(17, 10, 0.05536683332429467),
  @(17, 13, 0.7234496666666667),
  @(17, 14, 2.170349),
  @(17, 15, 0.1618662157416207),
  @(17, 16, 0.1618662157416207),
  @(17, 17, 0.138854588322),
  @(17, 18, 1.249691294898),
  @(17, 19, 0.008962019787800637),
  @(17, 20, 0.008962019787800637)
)

foreach ($item in $data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  $ws.Cells.Item($r, $c).Value = $v
}
